# Auto-generated Excel COM-interop script to apply cell value updates
# as described by the diff (odds data updates for Jogos_da_Semana_FlashScore_2025-04-15.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 1.18
$ws.Range("K2").Value = 4.5
$ws.Range("P2").Value = 1.85
$ws.Range("Q2").Value = 1.95
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
$ws.Range("N3").Value = 2.35
$ws.Range("O3").Value = 1.57
$ws.Range("J4").Value = 1.08
$ws.Range("K4").Value = 8
$ws.Range("J8").Value = 1.1
$ws.Range("K8").Value = 7
$ws.Range("G9").Value = 2.45
$ws.Range("I9").Value = 3.25
$ws.Range("U9").Value = 10
$ws.Range("AC9").Value = 81
$ws.Range("AF9").Value = 15
$ws.Range("AJ9").Value = 51
$ws.Range("G10").Value = 1.91
$ws.Range("H10").Value = 3.1
$ws.Range("J10").Value = 1.11
$ws.Range("K10").Value = 6.5
$ws.Range("L10").Value = 1.44
$ws.Range("M10").Value = 2.63
$ws.Range("N10").Value = 2.5
$ws.Range("O10").Value = 1.5
$ws.Range("P10").Value = 1.53
$ws.Range("Q10").Value = 2.38
$ws.Range("R10").Value = 2.2
$ws.Range("S10").Value = 1.62
$ws.Range("T10").Value = 5.5
$ws.Range("V10").Value = 9.5
$ws.Range("X10").Value = 19
$ws.Range("Y10").Value = 41
$ws.Range("Z10").Value = 6.5
$ws.Range("AA10").Value = 6
$ws.Range("AE10").Value = 9.5
$ws.Range("AG10").Value = 17
$ws.Range("N11").Value = 2.3
$ws.Range("O11").Value = 1.6
$ws.Range("G12").Value = 1.48
$ws.Range("H12").Value = 3.9
$ws.Range("I12").Value = 6
$ws.Range("U12").Value = 7
$ws.Range("AC12").Value = 51
$ws.Range("AE12").Value = 15
$ws.Range("AG12").Value = 19
$ws.Range("P13").Value = 1.47
$ws.Range("AD13").Value = 800
$ws.Range("P14").Value = 1.41
$ws.Range("Q14").Value = 2.62
$ws.Range("G15").Value = 2.9
$ws.Range("H15").Value = 3.4
$ws.Range("J15").Value = 1.04
$ws.Range("K15").Value = 13
$ws.Range("L15").Value = 1.25
$ws.Range("M15").Value = 3.75
$ws.Range("N15").Value = 1.83
$ws.Range("O15").Value = 2.03
$ws.Range("P15").Value = 1.33
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 1.67
$ws.Range("S15").Value = 2.1
$ws.Range("T15").Value = 11
$ws.Range("W15").Value = 29
$ws.Range("Y15").Value = 29
$ws.Range("Z15").Value = 12
$ws.Range("AD15").Value = 151
$ws.Range("AE15").Value = 9
$ws.Range("AF15").Value = 12
$ws.Range("AI15").Value = 17
$ws.Range("H17").Value = 3.15
$ws.Range("I17").Value = 3
$ws.Range("N17").Value = 2
$ws.Range("O17").Value = 1.65
$ws.Range("T17").Value = 6.1
$ws.Range("U17").Value = 8.75
$ws.Range("V17").Value = 7.6
$ws.Range("W17").Value = 17
$ws.Range("X17").Value = 15
$ws.Range("Y17").Value = 24
$ws.Range("Z17").Value = 8.5
$ws.Range("AB17").Value = 12
$ws.Range("AC17").Value = 55
$ws.Range("AD17").Value = 400
$ws.Range("AE17").Value = 7.2
$ws.Range("AF17").Value = 12.5
$ws.Range("AH17").Value = 30
$ws.Range("N18").Value = 2.05
$ws.Range("O18").Value = 1.75
$ws.Range("H19").Value = 4.5
$ws.Range("J19").Value = 1.05
$ws.Range("K19").Value = 11
$ws.Range("N19").Value = 1.88
$ws.Range("O19").Value = 1.98
$ws.Range("P19").Value = 1.36
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 2.05
$ws.Range("S19").Value = 1.7
$ws.Range("T19").Value = 6.5
$ws.Range("U19").Value = 6.5
$ws.Range("W19").Value = 9.5
$ws.Range("X19").Value = 13
$ws.Range("Y19").Value = 29
$ws.Range("Z19").Value = 11
$ws.Range("AB19").Value = 21
$ws.Range("AC19").Value = 67
$ws.Range("AD19").Value = 451
$ws.Range("AE19").Value = 15
$ws.Range("AG19").Value = 21
$ws.Range("AH19").Value = 81
$ws.Range("AJ19").Value = 51
$ws.Range("G20").Value = 2.7
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 2.5
$ws.Range("L20").Value = 1.3
$ws.Range("M20").Value = 3.4
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 1.8
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("R20").Value = 1.73
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 9
$ws.Range("Y20").Value = 29
$ws.Range("Z20").Value = 10
$ws.Range("AD20").Value = 201
$ws.Range("AE20").Value = 8.5
$ws.Range("AH20").Value = 23
$ws.Range("AI20").Value = 21
$ws.Range("AJ20").Value = 29
$ws.Range("G21").Value = 3.25
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 2.35
$ws.Range("N21").Value = 2.5
$ws.Range("O21").Value = 1.5
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 1.73
$ws.Range("W21").Value = 34
$ws.Range("X21").Value = 29
$ws.Range("Z21").Value = 6.5
$ws.Range("AE21").Value = 6.5
$ws.Range("AF21").Value = 10
$ws.Range("AH21").Value = 23
$ws.Range("AI21").Value = 23
